$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column C width: widened from a bestFit ~60.83 chars to a fixed ~68.83 chars ---
# (Do this BEFORE editing cell text so the page/print relationship isn't disturbed.)
$ws.Columns.Item(3).ColumnWidth = 68

# --- Update the embedded SQL queries: the joins moved from generic ".id" columns
#     to explicit "study_id" / "participant_id" columns on both sides of each join. ---
function Update-Query([string]$cellRef) {
    $t = $ws.Range($cellRef).Text
    $t2 = $t.Replace('prt ON std.id = prt."study.id"', 'prt ON std.study_id = prt."study.study_id"').Replace('dgn ON prt.id = dgn."participant.id"', 'dgn ON prt.participant_id = dgn."participant.participant_id"').Replace('trt ON prt.id = trt."participant.id"', 'trt ON prt.participant_id = trt."participant.participant_id"').Replace('trr ON prt.id = trr."participant.id"', 'trr ON prt.participant_id = trr."participant.participant_id"').Replace('srv ON prt.id = srv."participant.id"', 'srv ON prt.participant_id = srv."participant.participant_id"').Replace('rfs ON std.id = rfs."study.id"', 'rfs ON std.study_id = rfs."study.study_id"')
    $ws.Range($cellRef).Value = $t2
}

$cells = @("B2", "C2", "B3", "B4", "B5", "B6", "B7")
foreach ($c in $cells) {
    Update-Query $c
}
